$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.296676397323608
$ws.Range("B1").Value = 2.944759607315063
$ws.Range("C1").Value = 5.167603492736816
$ws.Range("D1").Value = 1.838336110115051
$ws.Range("E1").Value = 1.009628176689148
